$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "64.246.22"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -0.58%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.394.50"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -1.02%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "569.19"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.59%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "155.95"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -2.49%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.635"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +10.05%  "
$ws.Range("E8").Value = "  +0.05%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "3.401.66"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.90%  "
$ws.Range("E10").Value = "  -1.86%  "
$ws.Range("E11").Value = "  -1.76%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.439"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  +0.67%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "3.979.77"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -1.01%  "
$ws.Range("E14").Value = "  -0.11%  "
$ws.Range("E15").Value = "  -3.00%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "27.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -2.21%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.250.78"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -0.58%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "3.429.47"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +0.02%  "
$ws.Range("E19").Value = "  -0.79%  "
$ws.Range("E20").Value = "  -2.79%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "376.99"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.96%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "7.99"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -1.78%  "
$ws.Range("E23").Value = "  +0.22%  "
$ws.Range("E24").Value = "  -0.17%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "71.56"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -1.61%  "
$ws.Range("E26").Value = "  -3.79%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "10.30"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +8.18%  "
$ws.Range("E28").Value = "  -0.97%  "
$ws.Range("E29").Value = "  -0.28%  "
$ws.Range("E30").Value = "  +3.55%  "
$ws.Range("E31").Value = "  -0.79%  "
$ws.Range("E32").Value = "  -1.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "22.96"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.27%  "
$ws.Range("E34").Value = "  +0.71%  "
$ws.Range("E35").Value = "  +6.98%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "160.10"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -1.09%  "
$ws.Range("E37").Value = "  -0.96%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "6.92"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  +6.30%  "
$ws.Range("E39").Value = "  -0.74%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.869.03"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -5.06%  "
$ws.Range("E41").Value = "  +1.88%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "26.23"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -3.10%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "42.94"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.40%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.0312"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -1.10%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "25.83"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +5.59%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.763"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.51%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "321.50"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +6.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "0.111"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +4.21%  "
$ws.Range("E49").Value = "  -1.55%  "
$ws.Range("E50").Value = "  +0.72%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "6.52"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.14%  "
